# issue #5: stock data from json to db
#
# The stock ("股票") export now carries the extra bookkeeping columns that
# the json->db pipeline writes for every property-type sheet: a "category"
# field (normal/invalid filing bucket) right after property_category, and
# "source_file" / "index" provenance columns at the end of the row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column before the existing "date" column (column I) so
# "category" lands right after "property_category", pushing date /
# legislator_name / legislator_id one column to the right (I->J, J->K, K->L).
$ws.Columns.Item(9).Insert()

# Match the header/body formatting already used by the sheet (bold+border
# for row 1, plain body style for the data rows) instead of leaving the new
# cells with Excel's blank default.
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("H2:H4").Copy()
$ws.Range("I2:I4").PasteSpecial(-4122)

$ws.Cells.Item(1, 9).Value = "category"
$ws.Cells.Item(2, 9).Value = "normal"
$ws.Cells.Item(3, 9).Value = "normal"
$ws.Cells.Item(4, 9).Value = "normal"

# Append two more provenance columns: source_file and index.
$ws.Range("L1").Copy()
$ws.Range("M1:N1").PasteSpecial(-4122)
$ws.Range("L2:L4").Copy()
$ws.Range("M2:N4").PasteSpecial(-4122)

$ws.Cells.Item(1, 13).Value = "source_file"
$ws.Cells.Item(2, 13).Value = "tmp4fed1"
$ws.Cells.Item(3, 13).Value = "tmp4fed1"
$ws.Cells.Item(4, 13).Value = "tmp4fed1"

$ws.Cells.Item(1, 14).Value = "index"
$ws.Cells.Item(2, 14).Value = 69
$ws.Cells.Item(3, 14).Value = 70
$ws.Cells.Item(4, 14).Value = 71
